$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.89203699999999
$ws.Range("H2").Value = 107.676111
$ws.Range("I2").Value = 0.9301196142645664
$ws.Range("J2").Value = 0.9301196142645662
$ws.Range("M2").Value = 0.5273163333333333
$ws.Range("N2").Value = 1.581949
$ws.Range("Q2").Value = 18.92645734670433
$ws.Range("R2").Value = 170.338116120339
$ws.Range("S2").Value = 0.9301196142645664
$ws.Range("T2").Value = 0.9301196142645662

# Row 3
$ws.Range("I3").Value = 0.02407019339680195
$ws.Range("J3").Value = 0.02407019339680195
$ws.Range("M3").Value = 0.5273163333333333
$ws.Range("N3").Value = 1.581949
$ws.Range("Q3").Value = 0.4897902180158888
$ws.Range("R3").Value = 4.408111962142999
$ws.Range("S3").Value = 0.02407019339680195
$ws.Range("T3").Value = 0.02407019339680195

# Row 4
$ws.Range("G4").Value = 1.767752333333333
$ws.Range("H4").Value = 5.303257
$ws.Range("I4").Value = 0.04581019233863175
$ws.Range("J4").Value = 0.04581019233863175
$ws.Range("M4").Value = 0.5273163333333333
$ws.Range("N4").Value = 1.581949
$ws.Range("Q4").Value = 0.9321646786547777
$ws.Range("R4").Value = 8.389482107893
$ws.Range("S4").Value = 0.04581019233863175
$ws.Range("T4").Value = 0.04581019233863175
